$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1285457.5
$ws.Range("I17").Value = 150
$ws.Range("J17").Value = 1300578.8
$ws.Range("K17").Value = 450
$ws.Range("L17").Value = 3901736.4
$ws.Range("M17").Value = -282
$ws.Range("N17").Value = -3902072.4
$ws.Range("H40").Value = 1711.1111
$ws.Range("I40").Value = 1937.5
$ws.Range("J40").Value = 1530
$ws.Range("K40").Value = 1937.5
$ws.Range("L40").Value = 1530
$ws.Range("M40").Value = -1762.5
$ws.Range("N40").Value = -1880
$ws.Range("H96").Value = 3298.9546
$ws.Range("I96").Value = 4416.5
$ws.Range("J96").Value = 1957.9
$ws.Range("K96").Value = 13249.5
$ws.Range("L96").Value = 5873.700000000001
$ws.Range("M96").Value = -11876.5
$ws.Range("N96").Value = -8619.700000000001
$ws.Range("H129").Value = 2153.9583
$ws.Range("I129").Value = 430
$ws.Range("J129").Value = 2498.75
$ws.Range("K129").Value = 1290
$ws.Range("L129").Value = 7496.25
$ws.Range("M129").Value = 3710
$ws.Range("N129").Value = -17496.25
$ws.Range("H138").Value = 5954071
$ws.Range("I138").Value = 10102039
$ws.Range("J138").Value = 2637.6956
$ws.Range("K138").Value = 30306117
$ws.Range("L138").Value = 7913.0868
$ws.Range("M138").Value = -30300977
$ws.Range("N138").Value = -18193.0868
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11732.299
$ws.Range("I32").Value = 12525.6
$ws.Range("J32").Value = 9399.058999999999
$ws.Range("K32").Value = 12525.6
$ws.Range("L32").Value = 9399.058999999999
$ws.Range("M32").Value = -12238.6
$ws.Range("N32").Value = -9973.058999999999
$ws.Range("H122").Value = 5405.387
$ws.Range("I122").Value = 6382.0435
$ws.Range("J122").Value = 2597.5
$ws.Range("K122").Value = 19146.1305
$ws.Range("L122").Value = 7792.5
$ws.Range("M122").Value = -16696.1305
$ws.Range("N122").Value = -12692.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1277.6923
$ws.Range("I99").Value = 1268.3334
$ws.Range("J99").Value = 1285.7142
$ws.Range("K99").Value = 1268.3334
$ws.Range("L99").Value = 1285.7142
$ws.Range("M99").Value = 229.6666
$ws.Range("N99").Value = -4281.7142
$ws.Range("H134").Value = 2687.7646
$ws.Range("I134").Value = 1666.9286
$ws.Range("J134").Value = 3930.5217
$ws.Range("K134").Value = 5000.7858
$ws.Range("L134").Value = 11791.5651
$ws.Range("M134").Value = -2465.7858
$ws.Range("N134").Value = -16861.5651
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 4305604
$ws.Range("I6").Value = 4783782.5
$ws.Range("K6").Value = 4783782.5
$ws.Range("M6").Value = -4783669.5
$ws.Range("H7").Value = 125.625
$ws.Range("I7").Value = 100
$ws.Range("J7").Value = 141
$ws.Range("K7").Value = 100
$ws.Range("L7").Value = 141
$ws.Range("M7").Value = 13
$ws.Range("N7").Value = -367
$ws.Range("H31").Value = 5956410.5
$ws.Range("I31").Value = 4029.4666
$ws.Range("J31").Value = 30307060
$ws.Range("K31").Value = 4029.4666
$ws.Range("L31").Value = 30307060
$ws.Range("M31").Value = -3734.4666
$ws.Range("N31").Value = -30307650
$ws.Range("H34").Value = 5956410.5
$ws.Range("I34").Value = 4029.4666
$ws.Range("J34").Value = 30307060
$ws.Range("K34").Value = 4029.4666
$ws.Range("L34").Value = 30307060
$ws.Range("M34").Value = -3827.4666
$ws.Range("N34").Value = -30307464
$ws.Range("H41").Value = 15723.571
$ws.Range("I41").Value = 5000
$ws.Range("J41").Value = 17510.834
$ws.Range("K41").Value = 5000
$ws.Range("L41").Value = 17510.834
$ws.Range("M41").Value = -4572
$ws.Range("N41").Value = -18366.834
$ws.Range("H50").Value = 9513.200000000001
$ws.Range("I50").Value = 10000
$ws.Range("J50").Value = 9391.5
$ws.Range("K50").Value = 10000
$ws.Range("L50").Value = 9391.5
$ws.Range("M50").Value = -9375
$ws.Range("N50").Value = -10641.5
$ws.Range("H51").Value = 9225
$ws.Range("J51").Value = 9225
$ws.Range("L51").Value = 9225
$ws.Range("N51").Value = -10697
$ws.Range("H59").Value = 15588.333
$ws.Range("J59").Value = 15588.333
$ws.Range("L59").Value = 15588.333
$ws.Range("N59").Value = -17878.333
$ws.Range("H60").Value = 8041.1
$ws.Range("J60").Value = 8568.5
$ws.Range("L60").Value = 8568.5
$ws.Range("N60").Value = -9590.5
$ws.Range("H61").Value = 9225
$ws.Range("J61").Value = 9225
$ws.Range("L61").Value = 9225
$ws.Range("N61").Value = -9921
$ws.Range("H68").Value = 17958
$ws.Range("J68").Value = 17958
$ws.Range("L68").Value = 17958
$ws.Range("N68").Value = -19456
$ws.Range("H71").Value = 17958
$ws.Range("J71").Value = 17958
$ws.Range("L71").Value = 53874
$ws.Range("N71").Value = -61362
$ws.Range("H74").Value = 14157.286
$ws.Range("J74").Value = 14157.286
$ws.Range("L74").Value = 14157.286
$ws.Range("N74").Value = -15905.286
$ws.Range("H77").Value = 14157.286
$ws.Range("J77").Value = 14157.286
$ws.Range("L77").Value = 42471.858
$ws.Range("N77").Value = -51207.858
$ws.Range("H95").Value = 15905.75
$ws.Range("J95").Value = 15905.75
$ws.Range("L95").Value = 15905.75
$ws.Range("N95").Value = -21397.75
$ws.Range("H132").Value = 8622337
$ws.Range("I132").Value = 11112562
$ws.Range("K132").Value = 33337686
$ws.Range("M132").Value = -33335156
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 752.3043
$ws.Range("I113").Value = 436.32
$ws.Range("K113").Value = 1308.96
$ws.Range("M113").Value = 861.04
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3790.1282
$ws.Range("I132").Value = 2728.1785
$ws.Range("J132").Value = 6493.273
$ws.Range("K132").Value = 8184.5355
$ws.Range("L132").Value = 19479.819
$ws.Range("M132").Value = -5654.5355
$ws.Range("N132").Value = -24539.819
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1728.5714
$ws.Range("I68").Value = 1833.3334
$ws.Range("J68").Value = 1650
$ws.Range("K68").Value = 1833.3334
$ws.Range("L68").Value = 1650
$ws.Range("M68").Value = -1084.3334
$ws.Range("N68").Value = -3148
$ws.Range("H71").Value = 1728.5714
$ws.Range("I71").Value = 1833.3334
$ws.Range("J71").Value = 1650
$ws.Range("K71").Value = 9166.666999999999
$ws.Range("L71").Value = 8250
$ws.Range("M71").Value = -5422.666999999999
$ws.Range("N71").Value = -15738
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1162.5658
$ws.Range("I132").Value = 969.16394
$ws.Range("J132").Value = 1949.0667
$ws.Range("K132").Value = 2907.49182
$ws.Range("L132").Value = 5847.2001
$ws.Range("M132").Value = -377.4918200000002
$ws.Range("N132").Value = -10907.2001
